# Re-order the daily price records (rows 2-26) into their new positions.
# Each destination row takes on the full original row's data (columns A..R)
# from the source row indicated in the mapping below - i.e. the data rows
# were shuffled/re-sequenced (weekly logic) while headers stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: destination row -> source row (both refer to the ORIGINAL layout)
$rowMap = @{
    2  = 6
    3  = 15
    4  = 16
    5  = 2
    6  = 7
    7  = 5
    8  = 25
    9  = 19
    10 = 24
    11 = 13
    12 = 8
    13 = 20
    14 = 10
    15 = 11
    16 = 14
    17 = 9
    18 = 3
    19 = 22
    20 = 4
    21 = 23
    22 = 26
    23 = 18
    24 = 17
    25 = 12
    26 = 21
}

$firstCol = 1   # A
$lastCol  = 18  # R

# Snapshot every original row's values before we start overwriting, since
# the mapping is a permutation (sources and destinations overlap).
$snapshot = @{}
for ($r = 2; $r -le 26; $r++) {
    $rowVals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals += $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowVals
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c - $firstCol]
    }
}
